$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("scenarios")
$ws.Columns("D").Insert()
$ws.Range("D1").Value = "number_of_run"
$ws.Range("D1").Interior.ThemeColor = 1
$ws.Range("D1").Interior.Pattern = -4142
